# Revert "Powerpoint writer: consolidate text run nodes."
#
# Re-split the trailing-space-joined runs back into a separate word run
# and a separate single-space run, for each of the three title text
# frames that were affected:
#   Slide 1 title: "Header with inline code"
#       "Header " + "with "          -> "Header" + " " + "with" + " "
#   Slide 2 title: "Syntax highlighting"
#       "Syntax "                    -> "Syntax" + " "
#   Slide 3 title: "Two column slide"
#       "Two " + "column "           -> "Two" + " " + "column" + " "
#
# We do this by re-assigning (no-op, same character) a single-character
# Characters() sub-range at each space position: PowerPoint's COM
# TextRange model splits the run boundaries at the edited sub-range even
# when the replacement text equals the original text.

$p = $ppt.ActivePresentation

# --- Slide 1: "Header with inline code" ---
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(7, 1).Text = " "
$tr1.Characters(12, 1).Text = " "

# --- Slide 2: "Syntax highlighting" ---
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(7, 1).Text = " "

# --- Slide 3: "Two column slide" ---
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(4, 1).Text = " "
$tr3.Characters(11, 1).Text = " "
